# ----------------------------------------------------------------------------
# Adds the "2022-Q3" quarter to the 603987-康德莱 workbook:
#  1. Inserts a new summary row (row 2) into the "总计" sheet for 2022-Q3
#  2. Inserts a brand-new "2022-Q3" worksheet (placed right after "总计",
#     pushing every later quarter sheet one position to the right) and fills
#     it with the per-fund holdings detail for that quarter.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# xlCenter / xlTop alignment constants + continuous thin border style, reused below.
$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the 2022-Q3 row at row 2, above the
#    existing 2022-Q2 row, shifting everything else down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# The freshly-inserted row inherits stray formatting on B:D from the insert
# shift (it picks up the row-1 header look on some cells) - strip it so the
# new data cells come out unstyled, matching the rest of the data rows.
$summary.Range("B2:D2").ClearFormats()

# Copy the (still correctly-styled) bold/centered/bordered look from column A
# of the row that used to be row 2 (now pushed down to row 3) onto the new A2.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 29
$summary.Cells.Item(2, 4).Value = 9.95

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q3" worksheet with the per-fund holdings detail,
#    inserted immediately after "总计" (i.e. before the 2022-Q2 sheet).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q3.Name = "2022-Q3"
# Match the outline defaults ("summaryBelow"/"summaryRight") the other
# quarter sheets already carry in their <sheetPr><outlinePr .../></sheetPr>.
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop
    $cell.Borders.LineStyle = $xlContinuous
}

# Each inner array: @(rank, fund code, fund name, fund size, total stock
# position %, position share %, holding market value (billions), position rank)
$rows = @(
    @(0, "001645", "国泰大健康股票A", "36.92", "92.54", "6.87", "2.5364", 9),
    @(1, "020001", "国泰金鹰增长灵活配置混合", "20.78", "94.53", "7.73", "1.6063", 9),
    @(2, "009805", "国泰医药健康股票A", "8.71", "94.26", "9.85", "0.8579", 2),
    @(3, "000727", "融通健康产业灵活配置混合A", "22.64", "93.67", "3.63", "0.8218", 8),
    @(4, "009804", "国泰研究优势混合", "12.23", "92.87", "6.03", "0.7375", 10),
    @(5, "005760", "富国周期优势混合A", "22.68", "86.26", "2.84", "0.6441", 7),
    @(6, "009274", "融通健康产业灵活配置混合C", "17.64", "93.67", "3.63", "0.6403", 8),
    @(7, "160215", "国泰价值经典灵活配置混合（LOF）", "6.01", "93.77", "7.45", "0.4477", 8),
    @(8, "011335", "银河医药健康混合A", "7.55", "94.19", "4.60", "0.3473", 8),
    @(9, "008370", "国泰研究精选两年持有期混合", "2.98", "93.29", "8.29", "0.2470", 5),
    @(10, "002919", "东吴智慧医疗量化策略灵活配置混合A", "3.83", "90.19", "6.10", "0.2336", 5),
    @(11, "011321", "国泰大健康股票C", "3.20", "92.54", "6.87", "0.2198", 9),
    @(12, "011948", "东吴智慧医疗量化策略灵活配置混合C", "2.23", "90.19", "6.10", "0.1360", 5),
    @(13, "519673", "银河康乐股票A", "1.94", "93.82", "4.51", "0.0875", 7),
    @(14, "011876", "景顺长城医疗健康混合A", "2.10", "88.27", "4.09", "0.0859", 9),
    @(15, "002291", "诺安安鑫灵活配置混合", "2.66", "77.38", "2.68", "0.0713", 7),
    @(16, "013940", "东吴医疗服务股票A", "1.10", "91.62", "6.04", "0.0664", 5),
    @(17, "011326", "国泰医药健康股票C", "0.52", "94.26", "9.85", "0.0512", 2),
    @(18, "013941", "东吴医疗服务股票C", "0.51", "91.62", "6.04", "0.0308", 5),
    @(19, "011877", "景顺长城医疗健康混合C", "0.51", "88.27", "4.09", "0.0209", 9),
    @(20, "016018", "银河康乐股票C", "0.35", "93.82", "4.51", "0.0158", 7),
    @(21, "011565", "富国周期优势混合C", "0.41", "86.26", "2.84", "0.0116", 7),
    @(22, "013920", "兴华创新医疗6个月持有混合A", "0.18", "94.83", "5.41", "0.0097", 7),
    @(23, "014750", "兴华消费精选6个月持有混合A", "0.18", "94.78", "4.38", "0.0079", 6),
    @(24, "005210", "东吴双三角股票C", "0.10", "91.09", "6.09", "0.0061", 3),
    @(25, "005209", "东吴双三角股票A", "0.09", "91.09", "6.09", "0.0055", 3),
    @(26, "013921", "兴华创新医疗6个月持有混合C", "0.05", "94.83", "5.41", "0.0027", 7),
    @(27, "014751", "兴华消费精选6个月持有混合C", "0.02", "94.78", "4.38", "0.0009", 6),
    @(28, "015666", "银河医药健康混合C", "0.01", "94.19", "4.60", "0.0005", 8)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q3.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = $xlCenter
    $aCell.VerticalAlignment = $xlTop
    $aCell.Borders.LineStyle = $xlContinuous

    # Columns B-G (fund code/name/size/position/pct/value) are stored as text
    # in the source data, even though several look numeric - force text via
    # NumberFormat "@" before assignment, then clear the (now-stray) format
    # so the saved cell carries no style, matching the rest of the sheet.
    for ($c = 1; $c -le 5; $c++) {
        $cell = $q3.Cells.Item($r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$row[$c]
        $cell.ClearFormats()
    }
    $gCell = $q3.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = [string]$row[6]
    $gCell.ClearFormats()

    # Column H (position rank) is a genuine number.
    $hCell = $q3.Cells.Item($r, 8)
    $hCell.Value = $row[7]

    $r++
}
